# Adds a "2022-Q1" sheet (per-fund holdings) before the "总计" (totals)
# sheet, and records the new quarter's totals as the first data row of
# "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Scratch sheet used purely as a staging cell: writing "'<text>" to it and
# then Copy + PasteSpecial(xlPasteValues=-4163) into the real destination
# is the only reliable way to land a numeric-looking string (e.g. a fund
# code with leading zeros, or "94.29") as literal text instead of having
# Range.Value auto-coerce it into a number.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Add()
$scratch.Name = "__scratch__"
$stage = $scratch.Cells.Item(1, 1)

function Write-TextValue($targetCell, $text) {
    $stage.Value = "'" + $text
    $stage.Copy()
    $targetCell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header style (bold + border, same as the other per-fund sheets) lives on
# the "2021-Q4" sheet's header row; copy its format (not its text) over.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q1.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# Column-A style (bold + border, same as "A2" on the other sheets).
$q4.Range("A2").Copy()

$fundRows = @(
    @("870009", "广发资管平衡精选一年持有混合A", "11.34", "94.29", "5.43", "0.6158", 6),
    @("161609", "融通动力先锋混合",               "7.12",  "80.93", "2.96", "0.2108", 8),
    @("001990", "中欧数据挖掘多因子灵活配置混合A", "18.03", "84.18", "0.73", "0.1316", 4),
    @("872019", "广发资管平衡精选一年持有混合C", "1.54",  "94.29", "5.43", "0.0836", 6),
    @("001152", "融通新区域新经济灵活配置混合",   "2.81",  "80.98", "2.96", "0.0832", 9),
    @("002305", "光大保德信风格轮动混合A",         "3.99",  "90.94", "1.40", "0.0559", 10),
    @("007499", "光大保德信风格轮动混合C",         "2.82",  "90.94", "1.40", "0.0395", 10),
    @("004234", "中欧数据挖掘多因子灵活配置混合C", "5.06",  "84.18", "0.73", "0.0369", 4),
    @("006225", "人保量化基本面混合A",             "0.63",  "88.00", "1.50", "0.0094", 9),
    @("005166", "嘉实润和量化6个月定期开放混合",   "0.25",  "28.26", "0.62", "0.0016", 6),
    @("006157", "财通量化核心优选混合",             "0.09",  "92.85", "1.41", "0.0013", 9),
    @("006226", "人保量化基本面混合C",             "0.04",  "88.00", "1.50", "0.0006", 9)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = 2 + $i
    $row = $fundRows[$i]

    $q1.Cells.Item($r, 1).PasteSpecial(-4122)
    $q1.Cells.Item($r, 1).Value = $i

    Write-TextValue $q1.Cells.Item($r, 2) $row[0]
    Write-TextValue $q1.Cells.Item($r, 3) $row[1]
    Write-TextValue $q1.Cells.Item($r, 4) $row[2]
    Write-TextValue $q1.Cells.Item($r, 5) $row[3]
    Write-TextValue $q1.Cells.Item($r, 6) $row[4]
    Write-TextValue $q1.Cells.Item($r, 7) $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing 2021-Q4 / 2021-Q3 rows down by one.
# ---------------------------------------------------------------------
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A2:D3").ClearContents()

$summaryRows = @(
    @("2022-Q1", 12, 1.27),
    @("2021-Q4", 6, 2.28),
    @("2021-Q3", 2, 0.35)
)
for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = 2 + $i
    $row = $summaryRows[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# 3) Drop the scratch sheet used for text staging.
# ---------------------------------------------------------------------
$excel.DisplayAlerts = $false
$scratch.Delete()
